$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.417.98'
$ws.Range('E2').Value = '  -2.32%  '

$ws.Range('D3').Value = '1.649.83'
$ws.Range('E3').Value = '  -4.17%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.30%  '

$ws.Range('E6').Value = '  +0.19%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3618'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.43%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '47.50'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.05%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3272'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.33%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.119'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.51%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06928'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.43%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.00%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.924'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.34%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.21'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -8.26%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.587'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.03%  '

$ws.Range('D16').Value = '1.647.68'
$ws.Range('E16').Value = '  -4.26%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001037'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.95%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06505'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.30%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9994'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.00%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '76.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -9.45%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.906'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.70%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.66'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.50%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.62%  '

$ws.Range('D24').Value = '24.351.67'
$ws.Range('E24').Value = '  -2.57%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.426'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.80%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.338'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -16.74%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '146.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.58%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -11.57%  '

$ws.Range('D29').Value = '1.829.91'
$ws.Range('E29').Value = '  -4.25%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.15'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.17%  '

$ws.Range('E31').Value = '  -2.24%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.045'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.26%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.565'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -18.94%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08360'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.07%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.671'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.51%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.23'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -12.08%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.206'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.70%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06017'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.22%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02197'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.95%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.207'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.99%  '

$ws.Range('E41').Value = '  -7.71%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.157'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -9.43%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.13%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5795'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -10.40%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.724'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.91%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.54'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -9.91%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5528'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -10.14%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '121.61'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.60%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.932'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -10.11%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06887'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.66%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.64%  '
